# Update the "want to go" counts (column F) for the generated data output.
# Both the "展览" sheet and the "全部类型" sheet mirror the same rows and
# need identical updates:
#   F2: 252 -> 258
#   F3: 0   -> 4
#   F4: 249 -> 253

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 258
    $ws.Range("F3").Value = 4
    $ws.Range("F4").Value = 253
}
